# Refresh the crypto price/volume table (columns D and E) with the latest
# scraped figures, as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.019.44'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.460.53'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''527.20'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '''131.22'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('D7').Value = '''0.994'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = '''0.561'
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('D9').Value = '2.464.75'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').Value = '''0.0980'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('E11').Value = '  -3.08%  '
$ws.Range('D12').Value = '''4.97'
$ws.Range('E12').Value = '  -2.42%  '
$ws.Range('D13').Value = '''0.322'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = '2.892.57'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = '57.902.41'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '''21.87'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '2.453.46'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').Value = '''10.41'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').Value = '''4.16'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = '''317.47'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '''6.08'
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = '''64.89'
$ws.Range('E24').Value = '  +1.97%  '
$ws.Range('D25').Value = '''0.404'
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '''0.157'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').Value = '''7.29'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').Value = '''171.88'
$ws.Range('E29').Value = '  +2.87%  '
$ws.Range('D30').Value = '0.0₃0736'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').Value = '''1.16'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = '''0.998'
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('D36').Value = '''17.86'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('E37').Value = '  -4.04%  '
$ws.Range('D38').Value = '''3.82'
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('E39').Value = '  +1.68%  '
$ws.Range('D40').Value = '''36.20'
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').Value = '''0.803'
$ws.Range('E41').Value = '  +4.86%  '
$ws.Range('D42').Value = '''3.41'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = '''269.08'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').Value = '''126.90'
$ws.Range('E44').Value = '  +4.93%  '
$ws.Range('D45').Value = '''4.87'
$ws.Range('E45').Value = '  -3.20%  '
$ws.Range('D47').Value = '''0.0931'
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = '''0.0496'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').Value = '''0.0212'
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('D50').Value = '''16.49'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').Value = '1.716.55'
$ws.Range('E51').Value = '  -0.82%  '
